$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param(
        [string]$CellRef,
        [string]$Val
    )
    $rng = $ws.Range($CellRef)
    $origStyle = $rng.Style
    $rng.NumberFormat = "@"
    $rng.Value = $Val
    $rng.Style = $origStyle
}

Set-TextValue "D2" "26.123.58"
Set-TextValue "E2" "  -0.52%  "
Set-TextValue "D3" "1.656.45"
Set-TextValue "E3" "  -0.53%  "
Set-TextValue "D4" "1.002"
Set-TextValue "E4" "  -0.32%  "
Set-TextValue "D5" "218.30"
Set-TextValue "E5" "  -0.34%  "
Set-TextValue "D6" "0.5281"
Set-TextValue "E6" "  +0.98%  "
Set-TextValue "E7" "  -0.32%  "
Set-TextValue "D8" "0.2615"
Set-TextValue "E8" "  -1.90%  "
Set-TextValue "D9" "0.06342"
Set-TextValue "E9" "  +0.29%  "
Set-TextValue "D10" "20.48"
Set-TextValue "E10" "  -2.23%  "
Set-TextValue "D11" "0.07790"
Set-TextValue "E11" "  +0.53%  "
Set-TextValue "D12" "4.507"
Set-TextValue "E12" "  +1.44%  "
Set-TextValue "D13" "1.640.05"
Set-TextValue "E13" "  -1.84%  "
Set-TextValue "D14" "0.5486"
Set-TextValue "E14" "  +0.47%  "
Set-TextValue "E15" "  -0.64%  "
Set-TextValue "D16" "65.39"
Set-TextValue "E16" "  +0.83%  "
Set-TextValue "D17" "26.138.64"
Set-TextValue "E17" "  -0.56%  "
Set-TextValue "E18" "  -0.32%  "
Set-TextValue "D19" "4.562"
Set-TextValue "E19" "  -1.98%  "
Set-TextValue "D20" "192.98"
Set-TextValue "E20" "  -1.11%  "
Set-TextValue "E21" "  -0.61%  "
Set-TextValue "D22" "6.040"
Set-TextValue "E22" "  -0.45%  "
Set-TextValue "E23" "  -0.40%  "
Set-TextValue "D24" "141.16"
Set-TextValue "E24" "  +1.15%  "
Set-TextValue "D25" "0.1247"
Set-TextValue "E25" "  +0.69%  "
Set-TextValue "D26" "7.277"
Set-TextValue "E26" "  +1.12%  "
Set-TextValue "E27" "  +0.70%  "
Set-TextValue "E28" "  +1.26%  "
Set-TextValue "D29" "0.05920"
Set-TextValue "E29" "  -3.93%  "
Set-TextValue "D30" "1.280"
Set-TextValue "E30" "  -0.35%  "
Set-TextValue "D31" "3.528"
Set-TextValue "E31" "  -1.34%  "
Set-TextValue "D32" "3.245"
Set-TextValue "E32" "  -0.87%  "
Set-TextValue "D33" "1.570"
Set-TextValue "E33" "  -3.68%  "
Set-TextValue "D34" "0.9520"
Set-TextValue "E34" "  -2.32%  "
$ws.Range("B35").Value = "MXToken"
$ws.Range("C35").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
Set-TextValue "D35" "2.790"
Set-TextValue "E35" "  +0.09%  "
$ws.Range("B36").Value = "HuobiToken"
$ws.Range("C36").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
Set-TextValue "D36" "2.411"
Set-TextValue "D37" "0.5652"
Set-TextValue "E37" "  -1.35%  "
Set-TextValue "D38" "0.01614"
Set-TextValue "E38" "  +0.76%  "
Set-TextValue "D39" "5.817"
Set-TextValue "E39" "  -3.31%  "
Set-TextValue "D40" "0.8461"
Set-TextValue "E40" "  -1.41%  "
Set-TextValue "E41" "  -0.21%  "
Set-TextValue "D42" "101.73"
Set-TextValue "E42" "  +1.62%  "
Set-TextValue "D43" "1.016.68"
Set-TextValue "E43" "  -0.75%  "
Set-TextValue "D44" "1.800.37"
Set-TextValue "E44" "  -0.41%  "
Set-TextValue "D45" "57.14"
Set-TextValue "E45" "  -1.31%  "
Set-TextValue "E46" "  -1.78%  "
Set-TextValue "D47" "1.008"
Set-TextValue "E47" "  +0.12%  "
Set-TextValue "E48" "  +1.57%  "
Set-TextValue "D49" "1.479"
Set-TextValue "E49" "  -0.30%  "
Set-TextValue "D50" "0.05152"
Set-TextValue "E50" "  -0.69%  "
Set-TextValue "D51" "7.789"
Set-TextValue "E51" "  -3.65%  "
